$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the room-letter / HDJ labels in column A (rows 3-17).
# These cells keep their existing cell style; only the text content is removed.
$ws.Range("A3:A17").ClearContents()

# Keep the selection in sync with what the author left selected after the edit.
$ws.Range("A3:A17").Select()
